# Apply scheduled-runner price/profit updates to each class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1553.7368
$ws.Range("I15").Value = 1553.7368
$ws.Range("K15").Value = 4661.2104
$ws.Range("M15").Value = -4492.2104

# Row 62
$ws.Range("H62").Value = 8662.286
$ws.Range("I62").Value = 8007.091
$ws.Range("J62").Value = 11064.667
$ws.Range("K62").Value = 8007.091
$ws.Range("L62").Value = 11064.667
$ws.Range("M62").Value = -7383.091
$ws.Range("N62").Value = -12312.667

# Row 65
$ws.Range("H65").Value = 8662.286
$ws.Range("I65").Value = 8007.091
$ws.Range("J65").Value = 11064.667
$ws.Range("K65").Value = 40035.455
$ws.Range("L65").Value = 55323.335
$ws.Range("M65").Value = -36915.455
$ws.Range("N65").Value = -61563.335

# Row 88
$ws.Range("H88").Value = 6788.923
$ws.Range("J88").Value = 7568.727
$ws.Range("L88").Value = 7568.727
$ws.Range("N88").Value = -8380.726999999999

# Row 91
$ws.Range("H91").Value = 6788.923
$ws.Range("J91").Value = 7568.727
$ws.Range("L91").Value = 7568.727
$ws.Range("N91").Value = -10376.727

# Row 98
$ws.Range("H98").Value = 2322
$ws.Range("I98").Value = 2287.2942
$ws.Range("K98").Value = 2287.2942
$ws.Range("M98").Value = -789.2941999999998

# Row 118
$ws.Range("H118").Value = 948.93335
$ws.Range("I118").Value = 588.1429000000001
$ws.Range("K118").Value = 1764.4287
$ws.Range("M118").Value = -107.4287000000002

# Row 122
$ws.Range("H122").Value = 2322
$ws.Range("I122").Value = 2287.2942
$ws.Range("K122").Value = 6861.882599999999
$ws.Range("M122").Value = -4411.882599999999

# Row 132
$ws.Range("H132").Value = 21844.314
$ws.Range("I132").Value = 22369.146
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 67107.43799999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -64577.43799999999
$ws.Range("N132").Value = -17060

# Row 135
$ws.Range("H135").Value = 2586.65
$ws.Range("I135").Value = 2212.0833
$ws.Range("J135").Value = 3148.5
$ws.Range("K135").Value = 19908.7497
$ws.Range("L135").Value = 28336.5
$ws.Range("M135").Value = -17373.7497
$ws.Range("N135").Value = -33406.5

# Row 137
$ws.Range("H137").Value = 35469.668
$ws.Range("I137").Value = 75583
$ws.Range("K137").Value = 226749
$ws.Range("M137").Value = -224199

# Row 138
$ws.Range("H138").Value = 15996
$ws.Range("J138").Value = 38649.5
$ws.Range("L138").Value = 115948.5
$ws.Range("N138").Value = -126228.5


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 23475.412
$ws.Range("I32").Value = 23952.666
$ws.Range("J32").Value = 1999
$ws.Range("K32").Value = 23952.666
$ws.Range("L32").Value = 1999
$ws.Range("M32").Value = -23665.666
$ws.Range("N32").Value = -2573

# Row 97
$ws.Range("H97").Value = 1392.9584
$ws.Range("I97").Value = 1181.2
$ws.Range("J97").Value = 2451.75
$ws.Range("K97").Value = 1181.2
$ws.Range("L97").Value = 2451.75
$ws.Range("M97").Value = -685.2
$ws.Range("N97").Value = -3443.75

# Row 102
$ws.Range("H102").Value = 2516.5652
$ws.Range("I102").Value = 2529.3
$ws.Range("K102").Value = 2529.3
$ws.Range("M102").Value = -907.3000000000002

# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null

# Row 122
$ws.Range("H122").Value = 3455.4333
$ws.Range("I122").Value = 3348.3809
$ws.Range("K122").Value = 10045.1427
$ws.Range("M122").Value = -7595.1427

# Row 132
$ws.Range("H132").Value = 1141.5111
$ws.Range("I132").Value = 1002.9268
$ws.Range("K132").Value = 3008.7804
$ws.Range("M132").Value = -478.7803999999996


$ws = $wb.Worksheets.Item("BSM")
# Row 132
$ws.Range("H132").Value = 88223.8
$ws.Range("J132").Value = 88223.8
$ws.Range("L132").Value = 88223.8
$ws.Range("N132").Value = -98343.8

# Row 134
$ws.Range("H134").Value = 7863.93
$ws.Range("I134").Value = 9330.5
$ws.Range("J134").Value = 4479.5386
$ws.Range("K134").Value = 27991.5
$ws.Range("L134").Value = 13438.6158
$ws.Range("M134").Value = -25456.5
$ws.Range("N134").Value = -18508.6158


$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2274.375
$ws.Range("I16").Value = 1306.1818
$ws.Range("K16").Value = 1306.1818
$ws.Range("M16").Value = -1019.1818

# Row 31
$ws.Range("H31").Value = 12502275
$ws.Range("I31").Value = 14286598
$ws.Range("K31").Value = 14286598
$ws.Range("M31").Value = -14286303

# Row 34
$ws.Range("H34").Value = 12502275
$ws.Range("I34").Value = 14286598
$ws.Range("K34").Value = 14286598
$ws.Range("M34").Value = -14286396

# Row 58
$ws.Range("H58").Value = 12928.827
$ws.Range("I58").Value = 1301.6364
$ws.Range("K58").Value = 1301.6364
$ws.Range("M58").Value = -1098.6364

# Row 113
$ws.Range("H113").Value = 2274.375
$ws.Range("I113").Value = 1306.1818
$ws.Range("K113").Value = 1306.1818
$ws.Range("M113").Value = 863.8181999999999

# Row 132
$ws.Range("H132").Value = 73098.21000000001
$ws.Range("I132").Value = 101396.1
$ws.Range("J132").Value = 2353.5
$ws.Range("K132").Value = 304188.3
$ws.Range("L132").Value = 7060.5
$ws.Range("M132").Value = -301658.3
$ws.Range("N132").Value = -12120.5

# Row 134
$ws.Range("H134").Value = 1649.7354
$ws.Range("I134").Value = 1356.7667
$ws.Range("K134").Value = 4070.300099999999
$ws.Range("M134").Value = -1535.300099999999

# Row 136
$ws.Range("H136").Value = 12928.827
$ws.Range("I136").Value = 1301.6364
$ws.Range("K136").Value = 3904.9092
$ws.Range("M136").Value = -1354.9092


$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 86535240
$ws.Range("I4").Value = 115163496
$ws.Range("K4").Value = 345490488
$ws.Range("M4").Value = -345490376

# Row 122
$ws.Range("H122").Value = 1149.2858
$ws.Range("I122").Value = 870
$ws.Range("J122").Value = 1170.7693
$ws.Range("K122").Value = 7830
$ws.Range("L122").Value = 10536.9237
$ws.Range("M122").Value = -5380
$ws.Range("N122").Value = -15436.9237


$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2628.7144
$ws.Range("I97").Value = 1150.3334
$ws.Range("K97").Value = 1150.3334
$ws.Range("M97").Value = -654.3334

# Row 102
$ws.Range("H102").Value = 2436.077
$ws.Range("I102").Value = 1544.3334
$ws.Range("K102").Value = 1544.3334
$ws.Range("M102").Value = 77.66660000000002

# Row 126
$ws.Range("H126").Value = 3136.889
$ws.Range("I126").Value = 1872.1666
$ws.Range("K126").Value = 5616.4998
$ws.Range("M126").Value = -3146.4998

# Row 127
$ws.Range("H127").Value = 100000
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920

# Row 132
$ws.Range("H132").Value = 3462.25
$ws.Range("I132").Value = 3476.4167
$ws.Range("J132").Value = 3419.75
$ws.Range("K132").Value = 10429.2501
$ws.Range("L132").Value = 10259.25
$ws.Range("M132").Value = -7899.250100000001
$ws.Range("N132").Value = -15319.25


$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3110.8125
$ws.Range("I7").Value = 3277
$ws.Range("J7").Value = 2612.25
$ws.Range("K7").Value = 3277
$ws.Range("L7").Value = 2612.25
$ws.Range("M7").Value = -3165
$ws.Range("N7").Value = -2836.25

# Row 22
$ws.Range("H22").Value = 1298.4286
$ws.Range("I22").Value = 804
$ws.Range("K22").Value = 804
$ws.Range("M22").Value = -509

# Row 27
$ws.Range("H27").Value = 1298.4286
$ws.Range("I27").Value = 804
$ws.Range("K27").Value = 804
$ws.Range("M27").Value = -697

# Row 46
$ws.Range("H46").Value = 3463.1177
$ws.Range("I46").Value = 2200.25
$ws.Range("J46").Value = 3851.6924
$ws.Range("K46").Value = 2200.25
$ws.Range("L46").Value = 3851.6924
$ws.Range("M46").Value = -2012.25
$ws.Range("N46").Value = -4227.6924

# Row 116
$ws.Range("H116").Value = 99999.5
$ws.Range("J116").Value = 99999.5
$ws.Range("L116").Value = 99999.5
$ws.Range("N116").Value = -109177.5

# Row 126
$ws.Range("H126").Value = 3110.8125
$ws.Range("I126").Value = 3277
$ws.Range("J126").Value = 2612.25
$ws.Range("K126").Value = 9831
$ws.Range("L126").Value = 7836.75
$ws.Range("M126").Value = -7361
$ws.Range("N126").Value = -12776.75

# Row 132
$ws.Range("H132").Value = 4230.4
$ws.Range("I132").Value = 4236.75
$ws.Range("K132").Value = 12710.25
$ws.Range("M132").Value = -10180.25


$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 58772.04
$ws.Range("I122").Value = 90364.375
$ws.Range("J122").Value = 2607.889
$ws.Range("K122").Value = 271093.125
$ws.Range("L122").Value = 7823.667
$ws.Range("M122").Value = -268643.125
$ws.Range("N122").Value = -12723.667

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = $null
